$d = $word.ActiveDocument

# Locate the placeholder paragraph and remove it (text + paragraph mark)
$find = $d.Content
$find.Find.Execute('<Escopo do Projeto>', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$target = $find.Paragraphs(1)
$tr = $target.Range
$full = $d.Range($tr.Start, $tr.End + 1)
$full.Delete()

# Anchor: the paragraph that is now last before the trailing empty paragraph / sectPr
$anchor = $d.Paragraphs.Last

# paragraph 0: text paragraph with 5 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('Fundada em 1961, a WEG é de Jaraguá do Sul/SC, seu nome leva as iniciais dos seus 3 fundadores Werner Ricardo Voigt, Eggon João da Silva e Geraldo Werninghaus, que inicialmente produziam motores elétricos, começaram suas operações sob o Capital Social de Cr$ 3.600,00 (três mil e seiscentos cruzeiros). A partir da década de 80 ampliaram seu portifólio, incluindo a produção de componentes eletroeletrônicos, produtos para automação industrial, transformadores de força e distribuição, tintas líquidas e em pó e vernizes ')
$ip.Collapse(0)
$ip.InsertAfter('eletro isolantes')
$ip.Collapse(0)
$ip.InsertAfter('. Com essa expansão, a empresa se consolido')
$ip.Collapse(0)
$ip.InsertAfter('u')
$ip.Collapse(0)
$ip.InsertAfter(' como fornecedora de sistemas elétricos industriais completos. O primeiro parque fabril da empresa foi adquirido em 1964.')
$ip.Collapse(0)

# paragraph 1: empty paragraph with firstLine indent
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36

# paragraph 2: text paragraph with 3 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('Em 1970, começaram suas exportações, inicialmente para países como Guatemala, Uruguai, Paraguai Equador e Bolívia. No ano seguinte iniciaram as atividades de capital aberto na Bolsa de Valores. Ainda na década de 70, obteve feitos ')
$ip.Collapse(0)
$ip.InsertAfter('expressivos')
$ip.Collapse(0)
$ip.InsertAfter(' para sua expansão, construiu seu segundo parque fabril, atingiu a marca de 1 milhão de motores produzidos, expandiu as vendas para mais de 20 países e abriu um escritório na Alemanha.')
$ip.Collapse(0)

# paragraph 3: empty paragraph with firstLine indent
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36

# paragraph 4: text paragraph with 1 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('A década de 80 foi marcada pela expansão da área de atuação, tendo sido criadas WEG Acionamentos, WEG Transformadores, WEG Energia e WEG Química. Receberam também o primeiro prêmio “Revista Exame - Melhores e Maiores”. em 1989, o Sr. Eggon João da Silva sai da presidência e passa o cargo para Décio da Silva. Os fundadores formam o Conselho de Administração.')
$ip.Collapse(0)

# paragraph 5: empty paragraph with firstLine indent
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36

# paragraph 6: text paragraph with 7 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('Na ')
$ip.Collapse(0)
$ip.InsertAfter('década de')
$ip.Collapse(0)
$ip.InsertAfter(' 90, o que foi notório é a expansão para outros países, com a fundação das fabricas nos EUA, Alemanha, Inglaterra, França, Espanha e Suécia. Terminaram a década com um ')
$ip.Collapse(0)
$ip.InsertAfter('Market')
$ip.Collapse(0)
$ip.InsertAfter(' ')
$ip.Collapse(0)
$ip.InsertAfter('S')
$ip.Collapse(0)
$ip.InsertAfter('hare nacional de 79%, exportando 29% de sua produção para cerca de 55 países.')
$ip.Collapse(0)

# paragraph 7: empty paragraph with firstLine indent
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36

# paragraph 8: text paragraph with 3 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('De 2000, até os dias atuais, a empresa se consolidou como uma das maiores do mundo, em 2019 teve um faturamento de R$ 13,3 bi. Possui filiais em 36 ')
$ip.Collapse(0)
$ip.InsertAfter('países')
$ip.Collapse(0)
$ip.InsertAfter(' e fabricas em 12 países. O seu portifólio conta com mais de 600 produtos, com vendas presentes em 5 continentes. Conta com 31.800 colaboradores, mais de 1 bilhão de produtos de automação já fabricados e 16 milhões de motores produzidos anualmente.')
$ip.Collapse(0)

# paragraph 9: empty paragraph with firstLine indent
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36

# paragraph 10: text paragraph with 8 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('Visando atender as necessidades da ')
$ip.Collapse(0)
$ip.InsertAfter('indústria')
$ip.Collapse(0)
$ip.InsertAfter(' 4.0, a WEG está incluindo em seu ')
$ip.Collapse(0)
$ip.InsertAfter('portfólio')
$ip.Collapse(0)
$ip.InsertAfter(' produtos voltados para monitoramento total fabril através de IOT, com isso pensa em um produto para monitoramento de ')
$ip.Collapse(0)
$ip.InsertAfter('análise')
$ip.Collapse(0)
$ip.InsertAfter(' preditiva de motores')
$ip.Collapse(0)
$ip.InsertAfter('.')
$ip.Collapse(0)

# paragraph 11: bare empty paragraph
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last

# paragraph 12: text paragraph with 8 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('Para isso, será necessário desenvolver um')
$ip.Collapse(0)
$ip.InsertAfter(' sistema que faça o diagnóstico do funcionamento completo de motores do tamanho 63 a 450. O Sistema deverá fazer a medição da temperatura, análise de vibração,')
$ip.Collapse(0)
$ip.InsertAfter(' o tempo de funcionamento do motor, carga, velocidade, intervalo de lubrificação das engrenagens e mostrar os níveis de alerta de acordo com as manutenções preditivas que deverão ser realizadas posteriormente.')
$ip.Collapse(0)
$ip.InsertAfter(' ')
$ip.Collapse(0)
$ip.InsertAfter('Será')
$ip.Collapse(0)
$ip.InsertAfter(' disponibilizado para dispositivos inteligentes (disponível para os sistemas ')
$ip.Collapse(0)
$ip.InsertAfter('operacionais ')
$ip.Collapse(0)
$ip.InsertAfter('Android e iOS).')
$ip.Collapse(0)

# paragraph 13: text paragraph with 2 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('O sistema fará a coleta das informações descritas acima através de sensores e os dados serão enviados para a nuvem via Bluetooth ou Gatewa')
$ip.Collapse(0)
$ip.InsertAfter('y.')
$ip.Collapse(0)

# paragraph 14: text paragraph with 2 run(s)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$anchor.Format.FirstLineIndent = 36
$ip = $anchor.Range
$ip.Collapse(0)
$ip.InsertAfter('Após a coleta dos dados, as informações deverão ser enviadas para um servidor em nuvem, onde serão analisadas e transformadas em relatórios analíticos para a tomara rápida de decisões, que poderão ser analisados via smartphone, tablet ou desktop')
$ip.Collapse(0)
$ip.InsertAfter('.')
$ip.Collapse(0)

# paragraph 15: bare empty paragraph
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
